$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")
$ws.Activate()

# New cell G9: "Binary Exponentiation"
$ws.Range("C12").Value = "Binary S2 16"
$ws.Range("G9").Value = "Binary Exponentiation"
$ws.Range("G12").WrapText = $true
$ws.Range("G12").Value = "Google interview problem"
$ws.Range("D12").Value = "Stacks of coins"
$ws.Range("B12").Value = 8
$ws.Range("E12").Value = "No Link"

$ws.Rows.Item(12).RowHeight = 28.8

$ws.Range("G12").Select()
